# Update countries & provincias Spain
# Applies the 25-Jul-2020 14:08 -> 15:25 COVID data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Plain numeric refreshes - country stays on the same row, only its
#    statistics (Casos totales/Nuevos casos/Casos activos/Recuperados/
#    Casos criticos/Muertes hoy/Muertes) are updated.
# ---------------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4251024
$ws.Range("C4").Value = 2697
$ws.Range("E4").Value = 2074134
$ws.Range("G4").Value = 39
$ws.Range("H4").Value = 148529

# Row 6 - India
$ws.Range("B6").Value = 1358743
$ws.Range("C6").Value = 21721
$ws.Range("D6").Value = 866044
$ws.Range("E6").Value = 461059
$ws.Range("G6").Value = 234
$ws.Range("H6").Value = 31640

# Row 21 - Alemania
$ws.Range("B21").Value = 205983
$ws.Range("C21").Value = 23
$ws.Range("E21").Value = 6382

# Row 51 - Barein
$ws.Range("E51").Value = 3495
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 137

# Row 58 - Azerbaiyan
$ws.Range("B58").Value = 29633
$ws.Range("C58").Value = 321
$ws.Range("D58").Value = 22082
$ws.Range("E58").Value = 7143
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 408

# Row 62 - Serbia
$ws.Range("B62").Value = 23263
$ws.Range("C62").Value = 411
$ws.Range("E62").Value = 8690
$ws.Range("G62").Value = 8
$ws.Range("H62").Value = 526

# Row 65 - Uzbekistan
$ws.Range("B65").Value = 19755
$ws.Range("C65").Value = 395
$ws.Range("E65").Value = 9174
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 109

# Row 75 - Australia
$ws.Range("D75").Value = 9017
$ws.Range("E75").Value = 4786

# Row 82 - Republica de Macedonia
$ws.Range("B82").Value = 9934
$ws.Range("C82").Value = 137
$ws.Range("D82").Value = 5357
$ws.Range("E82").Value = 4117
$ws.Range("G82").Value = 9
$ws.Range("H82").Value = 460

# Row 98 - Croacia
$ws.Range("B98").Value = 4792
$ws.Range("C98").Value = 77
$ws.Range("D98").Value = 3778
$ws.Range("E98").Value = 881
$ws.Range("G98").Value = 5
$ws.Range("H98").Value = 133

# Row 114 - Sri Lanka
$ws.Range("B114").Value = 2768
$ws.Range("C114").Value = 4
$ws.Range("E114").Value = 654

# ---------------------------------------------------------------------------
# 2) Portugal overtakes Singapur (rows 45-46): Portugal gets refreshed
#    numbers and moves ahead of Singapur, whose own figures are unchanged.
# ---------------------------------------------------------------------------

$ws.Range("A45").Value = "Portugal"
$ws.Range("B45").Value = 49955
$ws.Range("C45").Value = 263
$ws.Range("D45").Value = 35010
$ws.Range("E45").Value = 13229
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 4
$ws.Range("H45").Value = 1716

$ws.Range("A46").Value = "Singapur"
$ws.Range("B46").Value = 49888
$ws.Range("C46").Value = 513
$ws.Range("D46").Value = 45172
$ws.Range("E46").Value = 4689
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 27

# ---------------------------------------------------------------------------
# 3) Zambia overtakes Paraguay and Grecia (rows 102-104): Zambia gets
#    refreshed numbers and moves to the top of the trio, pushing
#    Paraguay and Grecia down one row each (their own numbers unchanged).
# ---------------------------------------------------------------------------

$ws.Range("A102").Value = "Zambia"
$ws.Range("B102").Value = 4328
$ws.Range("C102").Value = 472
$ws.Range("D102").Value = 1953
$ws.Range("E102").Value = 2236
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 3
$ws.Range("H102").Value = 139

$ws.Range("A103").Value = "Paraguay"
$ws.Range("B103").Value = 4224
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 2596
$ws.Range("E103").Value = 1590
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 38

$ws.Range("A104").Value = "Grecia"
$ws.Range("B104").Value = 4135
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 1374
$ws.Range("E104").Value = 2560
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 201

# ---------------------------------------------------------------------------
# 4) Groenlandia and Islas Malvinas swap places (rows 210-211); both are
#    tied on every figure so only the labels need to change.
# ---------------------------------------------------------------------------

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Islas Malvinas"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# ---------------------------------------------------------------------------
# 5) Refresh the "last updated" timestamp banner in A1.
# ---------------------------------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 15:25"
